# Ads_Crawler.xlsx - "Updated for loading screens."
# Fills in the still-missing ads.txt lookup results (column C) for the
# rows that were added while the sheet was still loading, fixes one
# mis-recorded result, strips the leftover unused "Hyperlink" cell style
# from the few cells that still carried it, re-applies the AutoFilter on
# the table, and leaves the selection where the user ended up (G14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix a previously wrong result -----------------------------------
$ws.Range("C23").Value = "Yes"

# --- Strip the stray "Hyperlink" formatting that a few C-cells had ---
# (no actual hyperlinks on the sheet, just unused leftover styling)
$ws.Range("C21").Style = "Normal"
$ws.Range("C29").Style = "Normal"
$ws.Range("C45").Style = "Normal"
$ws.Range("C48").Style = "Normal"
$wb.Styles.Item("Hyperlink").Delete()

# --- Fill in the results for rows that finished loading ---------------
$results = @{
    48 = "No"
    49 = "No"
    50 = "No"
    51 = "Yes"
    52 = "Yes"
    53 = "Yes"
    54 = "Yes"
    55 = "Yes"
    56 = "Yes"
    57 = "Yes"
    58 = "Yes"
    59 = "Yes"
    60 = "Yes"
    61 = "Yes"
    62 = "Yes"
    63 = "Yes"
    64 = "No"
    65 = "No"
    66 = "No"
    67 = "No"
    68 = "Yes"
    69 = "Yes"
    70 = "Yes"
    71 = "No"
    72 = "Yes"
    73 = "Yes"
    74 = "Yes"
    75 = "Yes"
    76 = "No"
    77 = "No"
    78 = "No"
    79 = "Yes"
    80 = "No"
    81 = "No"
    82 = "Yes"
    83 = "Yes"
    84 = "No"
    85 = "Yes"
}

foreach ($row in $results.Keys | Sort-Object) {
    $ws.Range("C$row").Value = $results[$row]
}

# --- Re-apply the AutoFilter over the whole table ----------------------
$ws.Range("A1:F85").AutoFilter() | Out-Null

# --- Leave the selection where editing finished ------------------------
$ws.Range("G14").Select() | Out-Null
